# Row 2 ("A2:GG2") contains per-feature attribution scores that were
# recomputed (relative direction update). Replace the full row with the
# new values column-by-column (A=1 ... GG=189).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    -0.0,-0.0736143064681548,-0.0,0.2023283625086515,0.005361901149070607,-0.0,0.0,-0.0,-0.0,-0.0,0.005910374655943606,-0.0,
    0.2017405004068997,-0.003175673222564392,-0.0,0.0,0.0,-0.0,0.0,-0.0836551046379089,0.0,0.01416007321150566,-0.03298341659304817,0.0,
    -0.0,-0.0,-0.0,0.0,-0.05416417955287071,0.0,-0.01788162495550331,0.0002069792777307436,-0.0,-0.0,-0.0,0.0,
    -0.0,-0.03184932082569965,0.0,0.02827770634814052,0.06944358562979185,-0.0,0.0,-0.0,-0.0,0.0,-0.1497027310705481,0.0,
    0.07080831603100772,-0.001817919973327277,-0.0,0.0,0.0,0.0,-0.0,-0.01418805710578807,0.0,0.08649163433815991,0.03283123518905573,0.0,
    0.0,-0.0,-0.0,0.0,0.03190481457958391,0.0,-0.04257541708426302,-0.08905310676590357,0.0,0.0,0.0,0.0,
    0.0,-0.04640710802875297,0.0,0.01094127294829059,-0.02004983166574863,-0.0,-0.0,0.0,0.0,-0.0,0.03274941356648393,0.0,
    -0.03099117283538349,0.01593166186343111,0.0,-0.0,0.0,0.0,-0.0,-0.01061381960660221,-0.0,0.02133540246658532,0.03716094318380431,0.0,
    0.0,0.0,-0.0,-0.0,0.04597314766486385,-0.0,-0.03342648399499332,0.01017704690408558,-0.0,0.0,0.0,-0.0,
    -0.0,0.02874569132567836,-0.0,0.02900526664094873,0.03380050877759293,0.0,-0.0,-0.0,-0.0,0.0,-0.01950247745448723,-0.0,
    0.03637034262361485,-0.01945341551444906,-0.0,0.0,0.0,0.0,0.0,-0.05668836815106189,-0.0,-0.008248668484950638,-0.02434073422596091,0.0,
    -0.0,0.0,0.0,-0.0,0.04085971144248264,0.0,0.06686095049629477,-0.02407360759003618,0.0,0.0,0.0,0.0,
    0.0,0.04575396385905522,0.0,-0.0349068518066118,0.03654119765287879,0.0,-0.0,0.0,0.0,0.0,0.04340351386436194,0.0,
    -0.02669773779825179,0.01710822948871973,0.0,-0.0,0.0,0.0,-0.0,0.001691558334483329,0.0,-0.006288705109680439,-0.008476023910300627,-0.0,
    0.0,-0.0,-0.0,-0.0,-0.01406165917757605,-0.0,-0.0184191401413391,0.005010644615690384,0.0,-0.0,-0.0,-0.0,
    0.0,-0.03040407898268115,0.0,0.03399166782388464,-0.0,0.0,-0.0,-0.0,0.0
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $values[$i]
}
